$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "26.813.94"
Set-TextValue $ws.Range("E2") "  -1.37%  "
Set-TextValue $ws.Range("D3") "1.810.33"
Set-TextValue $ws.Range("E3") "  -0.15%  "
Set-TextValue $ws.Range("D4") "1.005"
Set-TextValue $ws.Range("E4") "  +0.28%  "
Set-TextValue $ws.Range("D5") "309.89"
Set-TextValue $ws.Range("E5") "  -0.75%  "
Set-TextValue $ws.Range("D6") "1.002"
Set-TextValue $ws.Range("E6") "  +0.02%  "
Set-TextValue $ws.Range("D7") "0.4640"
Set-TextValue $ws.Range("E7") "  +0.21%  "
Set-TextValue $ws.Range("D8") "0.3692"
Set-TextValue $ws.Range("E8") "  -1.34%  "
Set-TextValue $ws.Range("D9") "0.07356"
Set-TextValue $ws.Range("E9") "  -0.39%  "
Set-TextValue $ws.Range("D10") "0.8680"
Set-TextValue $ws.Range("E10") "  +0.19%  "
Set-TextValue $ws.Range("D11") "20.35"
Set-TextValue $ws.Range("E11") "  -1.04%  "
Set-TextValue $ws.Range("D12") "1.912.33"
Set-TextValue $ws.Range("E12") "  +5.51%  "
Set-TextValue $ws.Range("D13") "5.343"
Set-TextValue $ws.Range("E13") "  -0.67%  "
Set-TextValue $ws.Range("B14") "TRON"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D14") "0.07075"
Set-TextValue $ws.Range("E14") "  -0.12%  "
Set-TextValue $ws.Range("B15") "Litecoin"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D15") "91.59"
Set-TextValue $ws.Range("E15") "  -0.07%  "
Set-TextValue $ws.Range("B16") "Chainlink"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D16") "6.490"
Set-TextValue $ws.Range("E16") "  -2.40%  "
Set-TextValue $ws.Range("D17") "1.007"
Set-TextValue $ws.Range("E17") "  +0.35%  "
Set-TextValue $ws.Range("D18") "0.000008701"
Set-TextValue $ws.Range("E18") "  -0.40%  "
Set-TextValue $ws.Range("D19") "1.003"
Set-TextValue $ws.Range("E19") "  +0.13%  "
Set-TextValue $ws.Range("D20") "14.67"
Set-TextValue $ws.Range("E20") "  -1.32%  "
Set-TextValue $ws.Range("D21") "26.888.01"
Set-TextValue $ws.Range("E21") "  -1.09%  "
Set-TextValue $ws.Range("D22") "5.335"
Set-TextValue $ws.Range("E22") "  +0.60%  "
Set-TextValue $ws.Range("D23") "10.50"
Set-TextValue $ws.Range("E23") "  -3.64%  "
Set-TextValue $ws.Range("D24") "2.113.63"
Set-TextValue $ws.Range("E24") "  +3.13%  "
Set-TextValue $ws.Range("E25") "  -1.55%  "
Set-TextValue $ws.Range("D26") "151.63"
Set-TextValue $ws.Range("D27") "18.40"
Set-TextValue $ws.Range("E27") "  -0.65%  "
Set-TextValue $ws.Range("D28") "2.152"
Set-TextValue $ws.Range("D29") "5.279"
Set-TextValue $ws.Range("E29") "  +0.32%  "
Set-TextValue $ws.Range("D30") "115.35"
Set-TextValue $ws.Range("E30") "  -1.19%  "
Set-TextValue $ws.Range("D31") "0.08915"
Set-TextValue $ws.Range("E31") "  +0.36%  "
Set-TextValue $ws.Range("D32") "0.7556"
Set-TextValue $ws.Range("E32") "  -2.11%  "
Set-TextValue $ws.Range("D33") "1.154"
Set-TextValue $ws.Range("E33") "  -1.33%  "
Set-TextValue $ws.Range("B34") "HuobiToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D34") "2.932"
Set-TextValue $ws.Range("E34") "  +0.16%  "
Set-TextValue $ws.Range("B35") "Filecoin"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D35") "4.465"
Set-TextValue $ws.Range("E35") "  -0.99%  "
Set-TextValue $ws.Range("D36") "1.002"
Set-TextValue $ws.Range("E36") "  +0.03%  "
Set-TextValue $ws.Range("D37") "1.093"
Set-TextValue $ws.Range("E37") "  -1.67%  "
Set-TextValue $ws.Range("D38") "0.01951"
Set-TextValue $ws.Range("E38") "  -0.35%  "
Set-TextValue $ws.Range("D39") "0.05248"
Set-TextValue $ws.Range("E39") "  +0.26%  "
Set-TextValue $ws.Range("D40") "2.933"
Set-TextValue $ws.Range("E40") "  +0.90%  "
Set-TextValue $ws.Range("D41") "0.5338"
Set-TextValue $ws.Range("E41") "  +1.33%  "
Set-TextValue $ws.Range("D42") "7.198"
Set-TextValue $ws.Range("E42") "  -0.57%  "
Set-TextValue $ws.Range("D43") "2.362"
Set-TextValue $ws.Range("E43") "  -0.14%  "
Set-TextValue $ws.Range("D44") "0.1658"
Set-TextValue $ws.Range("E44") "  -1.27%  "
Set-TextValue $ws.Range("D45") "8.444"
Set-TextValue $ws.Range("E45") "  -1.45%  "
Set-TextValue $ws.Range("D46") "0.4940"
Set-TextValue $ws.Range("E46") "  -1.53%  "
Set-TextValue $ws.Range("D47") "10.31"
Set-TextValue $ws.Range("B48") "PaxDollar"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D48") "1.002"
Set-TextValue $ws.Range("E48") "  +0.00%  "
Set-TextValue $ws.Range("B49") "NEARProtocol"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D49") "1.672"
Set-TextValue $ws.Range("E49") "  +0.36%  "
Set-TextValue $ws.Range("D50") "103.02"
Set-TextValue $ws.Range("E50") "  -2.08%  "
Set-TextValue $ws.Range("D51") "0.06286"
Set-TextValue $ws.Range("E51") "  -0.54%  "
